$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.687.90"
$ws.Range("E2").Value = "  +5.14%  "

$ws.Range("D3").Value = "3.100.01"
$ws.Range("E3").Value = "  +3.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.20"
$ws.Range("E5").Value = "  +2.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.25"
$ws.Range("E6").Value = "  +10.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.092.88"
$ws.Range("E8").Value = "  +3.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.03"
$ws.Range("E10").Value = "  +17.21%  "

$ws.Range("E11").Value = "  +5.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +4.26%  "

$ws.Range("E13").Value = "  +4.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.28"
$ws.Range("E14").Value = "  +4.21%  "

$ws.Range("D15").Value = "3.605.30"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").Value = "64.628.13"
$ws.Range("E16").Value = "  +4.87%  "

$ws.Range("D17").Value = "3.101.65"
$ws.Range("E17").Value = "  +3.44%  "

$ws.Range("E18").Value = "  -0.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.82"
$ws.Range("E19").Value = "  +3.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.54"
$ws.Range("E20").Value = "  +1.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.83"
$ws.Range("E21").Value = "  +5.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.66"
$ws.Range("E22").Value = "  +10.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.675"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("E24").Value = "  +10.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.92"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.80"
$ws.Range("E27").Value = "  +4.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  +6.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.07"
$ws.Range("E29").Value = "  +8.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.11"
$ws.Range("E31").Value = "  +2.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  +2.85%  "

$ws.Range("E33").Value = "  +6.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  +3.88%  "

$ws.Range("E35").Value = "  +6.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.00"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "463.49"
$ws.Range("E37").Value = "  +4.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0408"
$ws.Range("E38").Value = "  +7.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0827"
$ws.Range("E39").Value = "  +4.75%  "

$ws.Range("D40").Value = "3.015.36"
$ws.Range("E40").Value = "  -3.60%  "

$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("E42").Value = "  +2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +15.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.72"
$ws.Range("E44").Value = "  +12.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  +8.01%  "

$ws.Range("E47").Value = "  +8.69%  "

$ws.Range("E48").Value = "  +4.43%  "

$ws.Range("E49").Value = "  +4.20%  "

$ws.Range("D50").Value = "0.0₃0515"
$ws.Range("E50").Value = "  +7.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.07"
$ws.Range("E51").Value = "  +3.24%  "

